# Event And Nurture Code
# Applies the data updates described in the commit to Sheet1:
#  - Emails row (row 3): admin/sent counts 1 -> 2
#  - Leads row (row 17): totals refreshed (482/0/0/0 -> 723/241/241/241)
#  - Tags row (row 20): value "3" -> "5"
#  - Two new summary rows appended: Event Programs (2), Nurture campaigns (3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Emails: Change Data Value / VALUE columns go from 1 to 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2

# Row 17 - Leads: refreshed totals
$ws.Range("B17").Value = 723
$ws.Range("C17").Value = 241
$ws.Range("D17").Value = 241
$ws.Range("E17").Value = 241

# Row 20 - Tags count, stored as text in the sheet ("3" -> "5")
$tags = $ws.Range("B20")
$tags.NumberFormat = "@"
$tags.Value = "5"
$tags.ClearFormats()

# New rows 28 & 29 - Event Programs / Nurture campaigns counts
$ws.Range("A28").Value = "Event Programs"
$ws.Range("B28").Value = 2
$ws.Range("A29").Value = "Nurture campaigns"
$ws.Range("B29").Value = 3

# Update sheet view: selection now targets the freshly appended rows and the
# sheet is no longer scrolled to A9 (it snaps back to the top-left cell A1)
$null = $ws.Range("A31:XFD32").Select()
